# Update the "dSF" column (F) values for the rows identified by the diff.
# Column A holds a 0-based index; spreadsheet row number = index + 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    4  = 2
    8  = -2
    12 = -3
    17 = 1
    19 = -2
    24 = -2
    30 = 6
    34 = 0
    42 = 0
    55 = -4
    57 = -2
    59 = 4
    60 = 6
    71 = 7
    76 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
